$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.017651052793478
$ws.Cells.Item(2, 4).Value = 1.050412519151949
$ws.Cells.Item(2, 5).Value = 1.018967393367267
$ws.Cells.Item(2, 6).Value = 1.053235327011146
$ws.Cells.Item(2, 9).Value = 1.0391899586383
$ws.Cells.Item(2, 10).Value = 1.022863957725663
$ws.Cells.Item(2, 11).Value = 1.053166867595929
$ws.Cells.Item(2, 12).Value = 1.021812037319454
$ws.Cells.Item(2, 13).Value = 1.055981855654603
$ws.Cells.Item(2, 14).Value = 1.011547502260672

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.019144097865797
$ws.Cells.Item(3, 4).Value = 1.051123604522355
$ws.Cells.Item(3, 5).Value = 1.020251949654212
$ws.Cells.Item(3, 6).Value = 1.054209255894751
$ws.Cells.Item(3, 9).Value = 1.039346313403227
$ws.Cells.Item(3, 10).Value = 1.023990130384912
$ws.Cells.Item(3, 11).Value = 1.053690655524277
$ws.Cells.Item(3, 12).Value = 1.022901200690109
$ws.Cells.Item(3, 13).Value = 1.056768376520259
$ws.Cells.Item(3, 14).Value = 1.011935495136237

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.020109184450622
$ws.Cells.Item(4, 4).Value = 1.051582789908887
$ws.Cells.Item(4, 5).Value = 1.021082501553459
$ws.Cells.Item(4, 6).Value = 1.054838529608173
$ws.Cells.Item(4, 9).Value = 1.039445795840487
$ws.Cells.Item(4, 10).Value = 1.024717506663779
$ws.Cells.Item(4, 11).Value = 1.054027975706567
$ws.Cells.Item(4, 12).Value = 1.023604785816203
$ws.Cells.Item(4, 13).Value = 1.057275768881446
$ws.Cells.Item(4, 14).Value = 1.012185687314277

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.020514671360561
$ws.Cells.Item(5, 4).Value = 1.051775607591244
$ws.Cells.Item(5, 5).Value = 1.021431516962199
$ws.Cells.Item(5, 6).Value = 1.055102856251585
$ws.Cells.Item(5, 9).Value = 1.039487213748191
$ws.Cells.Item(5, 10).Value = 1.025022981505681
$ws.Cells.Item(5, 11).Value = 1.054169400833127
$ws.Cells.Item(5, 12).Value = 1.023900295969128
$ws.Cells.Item(5, 13).Value = 1.057488709152601
$ws.Cells.Item(5, 14).Value = 1.01229066251524

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.020582740794117
$ws.Cells.Item(6, 4).Value = 1.051807969396715
$ws.Cells.Item(6, 5).Value = 1.021490109617275
$ws.Cells.Item(6, 6).Value = 1.055147224985467
$ws.Cells.Item(6, 9).Value = 1.039494144275629
$ws.Cells.Item(6, 10).Value = 1.025074253796473
$ws.Cells.Item(6, 11).Value = 1.054193124202926
$ws.Cells.Item(6, 12).Value = 1.023949897311929
$ws.Cells.Item(6, 13).Value = 1.057524441241256
$ws.Cells.Item(6, 14).Value = 1.012308276301725

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.020114603503061
$ws.Cells.Item(7, 4).Value = 1.051585367228046
$ws.Cells.Item(7, 5).Value = 1.021087165692105
$ws.Cells.Item(7, 6).Value = 1.054842062417042
$ws.Cells.Item(7, 9).Value = 1.039446350858096
$ws.Cells.Item(7, 10).Value = 1.024721589661882
$ws.Cells.Item(7, 11).Value = 1.054029866945893
$ws.Cells.Item(7, 12).Value = 1.023608735520207
$ws.Cells.Item(7, 13).Value = 1.057278615642505
$ws.Cells.Item(7, 14).Value = 1.012187090803995

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.01815584849633
$ws.Cells.Item(8, 4).Value = 1.050653027166143
$ws.Cells.Item(8, 5).Value = 1.019401651147646
$ws.Cells.Item(8, 6).Value = 1.053564662664255
$ws.Cells.Item(8, 9).Value = 1.039243149281764
$ws.Cells.Item(8, 10).Value = 1.02324483196542
$ws.Cells.Item(8, 11).Value = 1.053344216708644
$ws.Cells.Item(8, 12).Value = 1.022180371537683
$ws.Cells.Item(8, 13).Value = 1.05624798339046
$ws.Cells.Item(8, 14).Value = 1.011678806329567

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.014696190675046
$ws.Cells.Item(9, 4).Value = 1.049002974042252
$ws.Cells.Item(9, 5).Value = 1.016426415229034
$ws.Cells.Item(9, 6).Value = 1.051306601009455
$ws.Cells.Item(9, 9).Value = 1.038872142799558
$ws.Cells.Item(9, 10).Value = 1.020632168221781
$ws.Cells.Item(9, 11).Value = 1.052123709836175
$ws.Cells.Item(9, 12).Value = 1.019654195121091
$ws.Cells.Item(9, 13).Value = 1.054420042063081
$ws.Cells.Item(9, 14).Value = 1.010776449159257

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.012383870351809
$ws.Cells.Item(10, 4).Value = 1.047898134341725
$ws.Cells.Item(10, 5).Value = 1.014439145923316
$ws.Cells.Item(10, 6).Value = 1.049796368168493
$ws.Cells.Item(10, 9).Value = 1.038616104521948
$ws.Cells.Item(10, 10).Value = 1.018883070643067
$ws.Cells.Item(10, 11).Value = 1.051301754044857
$ws.Cells.Item(10, 12).Value = 1.017963589989448
$ws.Cells.Item(10, 13).Value = 1.053193387465066
$ws.Cells.Item(10, 14).Value = 1.010170278986373

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.011381114852248
$ws.Cells.Item(11, 4).Value = 1.04741858651621
$ws.Cells.Item(11, 5).Value = 1.013577668297082
$ws.Cells.Item(11, 6).Value = 1.049141251169549
$ws.Cells.Item(11, 9).Value = 1.03850317267662
$ws.Cells.Item(11, 10).Value = 1.018123885818309
$ws.Cells.Item(11, 11).Value = 1.050943869304699
$ws.Cells.Item(11, 12).Value = 1.017229935305094
$ws.Cells.Item(11, 13).Value = 1.052660311941882
$ws.Cells.Item(11, 14).Value = 1.009906687861509

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.011008411925949
$ws.Cells.Item(12, 4).Value = 1.047240288903112
$ws.Cells.Item(12, 5).Value = 1.013257524120045
$ws.Cells.Item(12, 6).Value = 1.048897733214643
$ws.Cells.Item(12, 9).Value = 1.038460914289505
$ws.Cells.Item(12, 10).Value = 1.017841612157929
$ws.Cells.Item(12, 11).Value = 1.050810638184165
$ws.Cells.Item(12, 12).Value = 1.016957175483936
$ws.Cells.Item(12, 13).Value = 1.052462013031343
$ws.Cells.Item(12, 14).Value = 1.009808608861123

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.011088368738128
$ws.Cells.Item(13, 4).Value = 1.047278542141905
$ws.Cells.Item(13, 5).Value = 1.01332620315867
$ws.Cells.Item(13, 6).Value = 1.048949976738828
$ws.Cells.Item(13, 9).Value = 1.03846999290883
$ws.Cells.Item(13, 10).Value = 1.017902173522154
$ws.Cells.Item(13, 11).Value = 1.050839230149782
$ws.Cells.Item(13, 12).Value = 1.017015694697061
$ws.Cells.Item(13, 13).Value = 1.052504562011264
$ws.Cells.Item(13, 14).Value = 1.009829654842074

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.011350311954206
$ws.Cells.Item(14, 4).Value = 1.04740385189384
$ws.Cells.Item(14, 5).Value = 1.013551208237666
$ws.Cells.Item(14, 6).Value = 1.049121125546831
$ws.Cells.Item(14, 9).Value = 1.038499685919529
$ws.Cells.Item(14, 10).Value = 1.018100558708424
$ws.Cells.Item(14, 11).Value = 1.050932862435901
$ws.Cells.Item(14, 12).Value = 1.017207393981431
$ws.Cells.Item(14, 13).Value = 1.052643926436049
$ws.Cells.Item(14, 14).Value = 1.009898584090169

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.011511672429045
$ws.Cells.Item(15, 4).Value = 1.047481036543854
$ws.Cells.Item(15, 5).Value = 1.013689820827011
$ws.Cells.Item(15, 6).Value = 1.049226552296349
$ws.Cells.Item(15, 9).Value = 1.038517939618567
$ws.Cells.Item(15, 10).Value = 1.01822275326769
$ws.Cells.Item(15, 11).Value = 1.050990513035749
$ws.Cells.Item(15, 12).Value = 1.017325473229727
$ws.Cells.Item(15, 13).Value = 1.052729754859872
$ws.Cells.Item(15, 14).Value = 1.009941031153553

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.012450387470904
$ws.Cells.Item(16, 4).Value = 1.047929936213691
$ws.Cells.Item(16, 5).Value = 1.014496298240564
$ws.Cells.Item(16, 6).Value = 1.049839821207071
$ws.Cells.Item(16, 9).Value = 1.038623555905295
$ws.Cells.Item(16, 10).Value = 1.018933416502375
$ws.Cells.Item(16, 11).Value = 1.051325464080678
$ws.Cells.Item(16, 12).Value = 1.018012245796871
$ws.Cells.Item(16, 13).Value = 1.053228725220302
$ws.Cells.Item(16, 14).Value = 1.010187748995121

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.013038809418879
$ws.Cells.Item(17, 4).Value = 1.048211212459247
$ws.Cells.Item(17, 5).Value = 1.015001914264137
$ws.Cells.Item(17, 6).Value = 1.050224192551607
$ws.Cells.Item(17, 9).Value = 1.03868925309947
$ws.Cells.Item(17, 10).Value = 1.019378706829916
$ws.Cells.Item(17, 11).Value = 1.051535041624098
$ws.Cells.Item(17, 12).Value = 1.018442604629894
$ws.Cells.Item(17, 13).Value = 1.053541199403931
$ws.Cells.Item(17, 14).Value = 1.010342208580642

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.01338188131379
$ws.Cells.Item(18, 4).Value = 1.048375165678792
$ws.Cells.Item(18, 5).Value = 1.015296737705221
$ws.Cells.Item(18, 6).Value = 1.050448276478565
$ws.Cells.Item(18, 9).Value = 1.038727373873138
$ws.Cells.Item(18, 10).Value = 1.019638262529798
$ws.Cells.Item(18, 11).Value = 1.051657094378446
$ws.Cells.Item(18, 12).Value = 1.018693470556472
$ws.Cells.Item(18, 13).Value = 1.053723274587201
$ws.Cells.Item(18, 14).Value = 1.010432194751889

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.01349883568724
$ws.Cells.Item(19, 4).Value = 1.048431050759282
$ws.Cells.Item(19, 5).Value = 1.01539724915276
$ws.Cells.Item(19, 6).Value = 1.050524664124461
$ws.Cells.Item(19, 9).Value = 1.038740338273051
$ws.Cells.Item(19, 10).Value = 1.019726734881196
$ws.Cells.Item(19, 11).Value = 1.051698678951648
$ws.Cells.Item(19, 12).Value = 1.018778983218298
$ws.Cells.Item(19, 13).Value = 1.053785326064474
$ws.Cells.Item(19, 14).Value = 1.010462859536628

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.01297569233979
$ws.Cells.Item(20, 4).Value = 1.048181045606488
$ws.Cells.Item(20, 5).Value = 1.014947676200743
$ws.Cells.Item(20, 6).Value = 1.050182964882818
$ws.Cells.Item(20, 9).Value = 1.038682225026677
$ws.Cells.Item(20, 10).Value = 1.019330949497237
$ws.Cells.Item(20, 11).Value = 1.051512575622158
$ws.Cells.Item(20, 12).Value = 1.018396447283995
$ws.Cells.Item(20, 13).Value = 1.053507693102779
$ws.Cells.Item(20, 14).Value = 1.010325647661355

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.01127318275996
$ws.Cells.Item(21, 4).Value = 1.047366956050272
$ws.Cells.Item(21, 5).Value = 1.013484954124866
$ws.Cells.Item(21, 6).Value = 1.049070731424638
$ws.Cells.Item(21, 9).Value = 1.038490950639662
$ws.Cells.Item(21, 10).Value = 1.018042146955803
$ws.Cells.Item(21, 11).Value = 1.050905298258043
$ws.Cells.Item(21, 12).Value = 1.017150950206742
$ws.Cells.Item(21, 13).Value = 1.052602895130223
$ws.Cells.Item(21, 14).Value = 1.009878290839852

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.010201384087388
$ws.Cells.Item(22, 4).Value = 1.046854109238528
$ws.Cells.Item(22, 5).Value = 1.012564394186409
$ws.Cells.Item(22, 6).Value = 1.048370393505714
$ws.Cells.Item(22, 9).Value = 1.038368892472737
$ws.Cells.Item(22, 10).Value = 1.017230210631401
$ws.Cells.Item(22, 11).Value = 1.050521761483553
$ws.Cells.Item(22, 12).Value = 1.016366420184726
$ws.Cells.Item(22, 13).Value = 1.052032328728962
$ws.Cells.Item(22, 14).Value = 1.009596038087464

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.01076969698784
$ws.Cells.Item(23, 4).Value = 1.047126073550252
$ws.Cells.Item(23, 5).Value = 1.013052486749304
$ws.Cells.Item(23, 6).Value = 1.048741754313265
$ws.Cells.Item(23, 9).Value = 1.038433768140426
$ws.Cells.Item(23, 10).Value = 1.017660788563112
$ws.Cells.Item(23, 11).Value = 1.050725244529458
$ws.Cells.Item(23, 12).Value = 1.016782452475188
$ws.Cells.Item(23, 13).Value = 1.052334956927773
$ws.Cells.Item(23, 14).Value = 1.009745759322595

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.013004212678254
$ws.Cells.Item(24, 4).Value = 1.04819467705047
$ws.Cells.Item(24, 5).Value = 1.014972184337161
$ws.Cells.Item(24, 6).Value = 1.050201594242519
$ws.Cells.Item(24, 9).Value = 1.038685401326194
$ws.Cells.Item(24, 10).Value = 1.019352529519974
$ws.Cells.Item(24, 11).Value = 1.05152272762969
$ws.Cells.Item(24, 12).Value = 1.018417304281278
$ws.Cells.Item(24, 13).Value = 1.053522833732275
$ws.Cells.Item(24, 14).Value = 1.01033313116018

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.015591598819808
$ws.Cells.Item(25, 4).Value = 1.049430399293832
$ws.Cells.Item(25, 5).Value = 1.017196227450191
$ws.Cells.Item(25, 6).Value = 1.051891215081676
$ws.Cells.Item(25, 9).Value = 1.038969590313986
$ws.Cells.Item(25, 10).Value = 1.021308872043597
$ws.Cells.Item(25, 11).Value = 1.052440699353794
$ws.Cells.Item(25, 12).Value = 1.020308394311141
$ws.Cells.Item(25, 13).Value = 1.054894018354295
$ws.Cells.Item(25, 14).Value = 1.011010533486326
